$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three rows whose Target cluster (column D) is "ECs":
# original rows 8 (MuSCs/ECs), 5 (FAPs/ECs), 2 (ECs/ECs) - delete from bottom to top
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Now rows are: 2=(ECs,FAPs) old3 ; 3=(ECs,MuSCs) old4 ; 4=(FAPs,FAPs) old6 ; 5=(FAPs,MuSCs) old7 ; 6=(MuSCs,FAPs) old9 ; 7=(MuSCs,MuSCs) old10

# Update the numeric values with the new TPM-derived data for each row
# Row 2: ECs -> FAPs
$ws.Range("G2").Value = 23.069913
$ws.Range("H2").Value = 69.20973899999998
$ws.Range("I2").Value = 0.9355059672894461
$ws.Range("J2").Value = 0.9355059672894461
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.183046666666667
$ws.Range("N2").Value = 3.54914
$ws.Range("O2").Value = 0.6222589862820888
$ws.Range("P2").Value = 0.6222589862820888
$ws.Range("Q2").Value = 27.29278367493999
$ws.Range("R2").Value = 245.6350530744599
$ws.Range("S2").Value = 0.5821269948663756
$ws.Range("T2").Value = 0.5821269948663756

# Row 3: ECs -> MuSCs
$ws.Range("G3").Value = 23.069913
$ws.Range("H3").Value = 69.20973899999998
$ws.Range("I3").Value = 0.9355059672894461
$ws.Range("J3").Value = 0.9355059672894461
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.718166
$ws.Range("N3").Value = 2.154498
$ws.Range("O3").Value = 0.3777410137179113
$ws.Range("P3").Value = 0.3777410137179112
$ws.Range("Q3").Value = 16.568027139558
$ws.Range("R3").Value = 149.112244256022
$ws.Range("S3").Value = 0.3533789724230705
$ws.Range("T3").Value = 0.3533789724230705

# Row 4: FAPs -> FAPs
$ws.Range("G4").Value = 0.7569533333333333
$ws.Range("H4").Value = 2.27086
$ws.Range("I4").Value = 0.03069514654402774
$ws.Range("J4").Value = 0.03069514654402774
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.183046666666667
$ws.Range("N4").Value = 3.54914
$ws.Range("O4").Value = 0.6222589862820888
$ws.Range("P4").Value = 0.6222589862820888
$ws.Range("Q4").Value = 0.895511117822222
$ws.Range("R4").Value = 8.059600060399999
$ws.Range("S4").Value = 0.01910033077226686
$ws.Range("T4").Value = 0.01910033077226686

# Row 5: FAPs -> MuSCs
$ws.Range("G5").Value = 0.7569533333333333
$ws.Range("H5").Value = 2.27086
$ws.Range("I5").Value = 0.03069514654402774
$ws.Range("J5").Value = 0.03069514654402774
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.718166
$ws.Range("N5").Value = 2.154498
$ws.Range("O5").Value = 0.3777410137179113
$ws.Range("P5").Value = 0.3777410137179112
$ws.Range("Q5").Value = 0.5436181475866666
$ws.Range("R5").Value = 4.89256332828
$ws.Range("S5").Value = 0.01159481577176088
$ws.Range("T5").Value = 0.01159481577176088

# Row 6: MuSCs -> FAPs
$ws.Range("G6").Value = 0.8334926666666668
$ws.Range("H6").Value = 2.500478
$ws.Range("I6").Value = 0.03379888616652608
$ws.Range("J6").Value = 0.03379888616652608
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.183046666666667
$ws.Range("N6").Value = 3.54914
$ws.Range("O6").Value = 0.6222589862820888
$ws.Range("P6").Value = 0.6222589862820888
$ws.Range("Q6").Value = 0.9860607209911112
$ws.Range("R6").Value = 8.87454648892
$ws.Range("S6").Value = 0.02103166064344623
$ws.Range("T6").Value = 0.02103166064344623

# Row 7: MuSCs -> MuSCs
$ws.Range("G7").Value = 0.8334926666666668
$ws.Range("H7").Value = 2.500478
$ws.Range("I7").Value = 0.03379888616652608
$ws.Range("J7").Value = 0.03379888616652608
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.718166
$ws.Range("N7").Value = 2.154498
$ws.Range("O7").Value = 0.3777410137179113
$ws.Range("P7").Value = 0.3777410137179112
$ws.Range("Q7").Value = 0.5985860944493334
$ws.Range("R7").Value = 5.387274850044
$ws.Range("S7").Value = 0.01276722552307985
$ws.Range("T7").Value = 0.01276722552307985

$wb.Save()
